$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column J: copy H3's current format (right-edge header border) to J3 first,
# --- before we repoint G3/H3 to I3's format.
$ws.Range("H3").Copy()
$ws.Range("J3").PasteSpecial(-4122)

# --- G3 and H3 now take on I3's header style (shared "inner" border).
$ws.Range("I3").Copy()
$ws.Range("G3:H3").PasteSpecial(-4122)

# --- J4 takes on I4's numeric/error cell style.
$ws.Range("I4").Copy()
$ws.Range("J4").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Values / formulas ---
# I3 header text changes.
$ws.Range("I3").Value = "Sum by FM change last green elements"
# New J3 header text.
$ws.Range("J3").Value = "Sum staticly using ;func=SUM"
# New J4 literal text value.
$ws.Range("J4").Value = "{R-T-SUM;func=SUM}"

# --- Column width: extend the I column's custom width onto the new J column. ---
$ws.Columns.Item(10).ColumnWidth = 18.67

# --- Selection marker moves from I1 to K1. ---
$ws.Range("K1").Select()
